$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 16 ("6.4.3 Identifying Subsystems") - Content Placeholder 2
#   * shrink the placeholder's height (ext cy 3777622 -> 2695998)
#   * remove the third bullet's text/level, leaving an empty paragraph
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$shContent = $s16.Shapes.Item(4)

# Resize the placeholder (EMU -> points, 1 pt = 12700 EMU; nudge by +0.5 EMU
# worth of points so the float32-precision COM round-trip lands on the exact
# target EMU value instead of one unit short).
$shContent.Height = (2695998 + 0.5) / 12700

# Clear the text of the third paragraph ("initial subsystem should be
# derived ...") while keeping the paragraph itself (now empty), and reset
# its outline level back to the top level.
$tr = $shContent.TextFrame.TextRange
$thirdPara = $tr.Characters(106, 120)
$thirdPara.IndentLevel = 1
$thirdPara.Text = ""

# ---------------------------------------------------------------------------
# Slide 17 - Picture 3
#   * nudge the picture down slightly (off y 584610 -> 595899)
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$shPic = $s17.Shapes.Item(1)
$shPic.Top = (595899 + 0.5) / 12700
